$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text (A1:D1), drop column E
$ws.Range("A1").Value = "ProductNumber"
$ws.Range("B1").Value = "ProductName"
$ws.Range("C1").Value = "QtyPerPackage"
$ws.Range("D1").Value = "Notes"
$ws.Range("E1").ClearContents()

# Column widths + base style for data entry columns A:D
$cols = $ws.Range("A1:D1").EntireColumn
$cols.ColumnWidth = 30
$cols.HorizontalAlignment = -4108
$cols.Locked = $false

# Header row formatting
$hdr = $ws.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.Font.Size = 14
$hdr.Interior.Color = 65535
$hdr.Borders.Weight = -4138
$hdr.HorizontalAlignment = -4108

# Protect the worksheet
$ws.Protect()

Write-Host "done"
